$wb = $excel.ActiveWorkbook

# --- Add 2022-Q1 sheet before the 总计 (Total) sheet ---
$wsRef = $wb.Worksheets.Item("2021-Q4")
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsQ1 = $wb.Worksheets.Add($wsLast)
$wsQ1.Name = "2022-Q1"

# copy header-row format (bold + border) and column-A format from an existing quarter sheet
$wsRef.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)
$wsRef.Range("A2").Copy()
$wsQ1.Range("A2:A23").PasteSpecial(-4122)

$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

$q1Data = @(
  @("002685", "中欧丰泓沪港深灵活配置混合A", "59.37", "94.70", "8.94", "5.3077", 1),
  @("005847", "富国沪港深业绩驱动混合A", "44.03", "74.38", "9.44", "4.1564", 1),
  @("501087", "交银施罗德瑞丰三年封闭运作混合型", "40.75", "92.66", "8.81", "3.5901", 1),
  @("010583", "富国蓝筹精选股票（QDII）美元", "14.71", "94.57", "9.55", "1.4048", 1),
  @("007455", "富国蓝筹精选股票（QDII）人民币", "14.71", "94.57", "9.55", "1.4048", 1),
  @("013991", "中欧港股通精选一年持有混合A", "12.87", "94.50", "5.76", "0.7413", 3),
  @("002686", "中欧丰泓沪港深灵活配置混合C", "7.65", "94.70", "8.94", "0.6839", 1),
  @("009984", "鹏华启航两年封闭运作混合", "14.73", "89.97", "4.26", "0.6275", 5),
  @("009846", "富兰克林国海港股通远见价值混合", "19.47", "86.72", "3.01", "0.5860", 5),
  @("011635", "富国港股通策略精选混合型证券投资基金A", "7.67", "76.88", "5.89", "0.4518", 2),
  @("013992", "中欧港股通精选一年持有混合C", "5.32", "94.50", "5.76", "0.3064", 3),
  @("011117", "富国沪港深业绩驱动混合C", "2.39", "74.38", "9.44", "0.2256", 1),
  @("862001", "光大阳光香港精选混合型集合资产管理计划（QDII）A 人民币", "3.91", "89.45", "4.02", "0.1572", 8),
  @("862011", "光大阳光香港精选混合型集合资产管理计划（QDII）A 美元", "3.91", "89.45", "4.02", "0.1572", 8),
  @("862012", "光大阳光香港精选混合型集合资产管理计划（QDII）C 人民币", "3.91", "89.45", "4.02", "0.1572", 8),
  @("011349", "淳厚现代服务业股票A", "3.58", "81.51", "2.53", "0.0906", 7),
  @("160125", "南方香港优选股票QDII-LOF", "2.46", "91.14", "3.51", "0.0863", 8),
  @("011636", "富国港股通策略精选混合型证券投资基金C", "0.92", "76.88", "5.89", "0.0542", 2),
  @("008134", "鹏华优选价值股票", "1.86", "92.62", "2.71", "0.0504", 9),
  @("004099", "前海开源沪港深景气行业精选灵活配置混合", "0.41", "93.07", "9.66", "0.0396", 1),
  @("011350", "淳厚现代服务业股票C", "0.63", "81.51", "2.53", "0.0159", 7),
  @("000761", "国富健康优质生活股票", "0.08", "84.48", "6.47", "0.0052", 2),
)

$r = 2
foreach ($row in $q1Data) {
    $wsQ1.Cells.Item($r, 1).Value = $r - 2
    $wsQ1.Cells.Item($r, 2).Value = "'" + $row[0]
    $wsQ1.Cells.Item($r, 3).Value = "'" + $row[1]
    $wsQ1.Cells.Item($r, 4).Value = "'" + $row[2]
    $wsQ1.Cells.Item($r, 5).Value = "'" + $row[3]
    $wsQ1.Cells.Item($r, 6).Value = "'" + $row[4]
    $wsQ1.Cells.Item($r, 7).Value = "'" + $row[5]
    $wsQ1.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Remove the stray "quote-prefix" formatting that typing a leading apostrophe
# applies, so the text cells end up with no explicit style (matches the
# original detail-sheet layout where only the header row and column A carry
# a style).
$wsQ1.Range("B2:G23").ClearFormats()

# --- Update 总计 (Total) sheet: insert a new 2022-Q1 row at the top of the data ---
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# the inserted row picks up formatting from the header row above; reset it,
# then re-apply column-A's normal style (same as the rows below it)
$wsTotal.Range("A2:D2").ClearFormats()
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 22
$wsTotal.Cells.Item(2, 4).Value = 20.3

for ($r = 3; $r -le 6; $r++) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
}
